# Add a new "SIDEBAR_SUBMENU_SUBMENU" column to the DPLKKPS016-001 sheet.
#
# The header row (row 1) gains its new column one slot to the right of where
# the data row (row 2) gains its new column, exactly mirroring the source
# edit: row 1's new cell lands in M1 (pushing old M1/N1 -> N1/O1), while row
# 2's new cell lands in L2 (pushing old L2/M2 -> M2/N2).
#
# NOTE: in this host, setting .Value on a range AFTER a PasteSpecial(formats)
# clobbers the just-pasted style back to the default, so every cell below
# sets its .Value FIRST and only THEN copies over the donor cell's format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values that need to slide right before they get overwritten.
$n1val = $ws.Range("N1").Value()
$m1val = $ws.Range("M1").Value()
$m2val = $ws.Range("M2").Value()
$l2val = $ws.Range("L2").Value()

# --- Row 1 (headers): shift M1/N1 -> N1/O1, then drop the new text into M1 ---
$ws.Range("O1").Value = $n1val
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("N1").Value = $m1val
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

$ws.Range("M1").Value = "SIDEBAR_SUBMENU_SUBMENU"
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

# --- Row 2 (data): shift L2/M2 -> M2/N2, then drop the new text into L2 ---
$ws.Range("N2").Value = $m2val
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

$ws.Range("M2").Value = $l2val
$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

$ws.Range("L2").Value = "Setup Kelengkapan Kepesertaan"
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column widths: the new column (L) gets a plain custom width; the
#     widths that used to belong to L/M/N/O slide one column to the right ---
$ws.Columns.Item(16).ColumnWidth = 21.666666666666668
$ws.Columns.Item(15).ColumnWidth = 16
$ws.Columns.Item(14).ColumnWidth = 15.166666666666666
$ws.Columns.Item(13).ColumnWidth = 17.5
$ws.Columns.Item(12).ColumnWidth = 14.166666666666666

# Update the active selection to match the saved view state
$ws.Range("L9").Select()
